# Update the "取得日時" (retrieved datetime) column to reflect the new
# append timestamp: 2026-02-07 02:24:07 (JST), replacing the previous
# 2026-02-07 01:49:50 value for every existing data row in the
# "ランサーズ" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-07 02:24:07"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
